# Converts an "RRGGBB" hex string into the BGR-packed integer that the
# PowerPoint COM RGB() convention expects (val = R + G*256 + B*65536).
function Hex2Rgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table on slide 5 switches to a different built-in table style.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{F61C2E9F-B32B-40B8-8E84-F06BAABE2A82}")

# ---------------------------------------------------------------------
# 2) The presentation's theme colour scheme switches from the custom
#    "Integral" (Red Violet) palette to the default Office theme palette.
# ---------------------------------------------------------------------
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

$newColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $newColors.Count; $i++) {
    $themeColors.Item($i).RGB = Hex2Rgb($newColors[$i - 1])
}
